$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3843.3635
$ws.Range("I64").Value = 2941.5715
$ws.Range("J64").Value = 5421.5
$ws.Range("K64").Value = 2941.5715
$ws.Range("L64").Value = 5421.5
$ws.Range("M64").Value = -2693.5715
$ws.Range("N64").Value = -5917.5
$ws.Range("H67").Value = 3843.3635
$ws.Range("I67").Value = 2941.5715
$ws.Range("J67").Value = 5421.5
$ws.Range("K67").Value = 2941.5715
$ws.Range("L67").Value = 5421.5
$ws.Range("M67").Value = -2083.5715
$ws.Range("N67").Value = -7137.5
$ws.Range("H116").Value = 13839619
$ws.Range("I116").Value = 19769784
$ws.Range("J116").Value = 2566.6667
$ws.Range("K116").Value = 19769784
$ws.Range("L116").Value = 2566.6667
$ws.Range("M116").Value = -19766342
$ws.Range("N116").Value = -9450.6667
$ws.Range("H125").Value = 14016514
$ws.Range("I125").Value = 2032
$ws.Range("J125").Value = 16018583
$ws.Range("K125").Value = 18288
$ws.Range("L125").Value = 144167247
$ws.Range("M125").Value = -15828
$ws.Range("N125").Value = -144172167
$ws.Range("H133").Value = 42257.418
$ws.Range("J133").Value = 42257.418
$ws.Range("L133").Value = 42257.418
$ws.Range("N133").Value = -52377.418
$ws.Range("H137").Value = 30303988
$ws.Range("I137").Value = 40000880
$ws.Range("J137").Value = 1198.25
$ws.Range("K137").Value = 120002640
$ws.Range("L137").Value = 3594.75
$ws.Range("M137").Value = -120000090
$ws.Range("N137").Value = -8694.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 7574.3
$ws.Range("I74").Value = 1199.0667
$ws.Range("J74").Value = 26700
$ws.Range("K74").Value = 1199.0667
$ws.Range("L74").Value = 26700
$ws.Range("M74").Value = -325.0667000000001
$ws.Range("N74").Value = -28448
$ws.Range("H77").Value = 7574.3
$ws.Range("I77").Value = 1199.0667
$ws.Range("J77").Value = 26700
$ws.Range("K77").Value = 5995.333500000001
$ws.Range("L77").Value = 133500
$ws.Range("M77").Value = -1627.333500000001
$ws.Range("N77").Value = -142236

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1835.091
$ws.Range("I94").Value = 1134
$ws.Range("K94").Value = 1134
$ws.Range("M94").Value = -683
$ws.Range("H134").Value = 27780276
$ws.Range("I134").Value = 35715704
$ws.Range("J134").Value = 6276.5
$ws.Range("K134").Value = 107147112
$ws.Range("L134").Value = 18829.5
$ws.Range("M134").Value = -107144577
$ws.Range("N134").Value = -23899.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H109").Value = 33750
$ws.Range("J109").Value = 33750
$ws.Range("L109").Value = 33750
$ws.Range("N109").Value = -35830

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 495
$ws.Range("I40").Value = 290
$ws.Range("J40").Value = 700
$ws.Range("K40").Value = 1160
$ws.Range("L40").Value = 2800
$ws.Range("M40").Value = -1091
$ws.Range("N40").Value = -2938
$ws.Range("H58").Value = 7206.5293
$ws.Range("J58").Value = 9153.846
$ws.Range("L58").Value = 27461.538
$ws.Range("N58").Value = -27717.538
$ws.Range("H64").Value = 3930.8572
$ws.Range("I64").Value = 1256
$ws.Range("K64").Value = 3768
$ws.Range("M64").Value = -3498
$ws.Range("H67").Value = 3930.8572
$ws.Range("I67").Value = 1256
$ws.Range("K67").Value = 3768
$ws.Range("M67").Value = -2832
$ws.Range("H75").Value = 1665.6
$ws.Range("J75").Value = 1985.909
$ws.Range("L75").Value = 5957.727000000001
$ws.Range("N75").Value = -7953.727000000001
$ws.Range("H78").Value = 1665.6
$ws.Range("J78").Value = 1985.909
$ws.Range("L78").Value = 17873.181
$ws.Range("N78").Value = -27857.181
$ws.Range("H121").Value = 676.8333
$ws.Range("I121").Value = 230
$ws.Range("J121").Value = 900.25
$ws.Range("K121").Value = 690
$ws.Range("L121").Value = 2700.75
$ws.Range("M121").Value = 620
$ws.Range("N121").Value = -5320.75
$ws.Range("H122").Value = 883.1667
$ws.Range("I122").Value = 300
$ws.Range("J122").Value = 1299.7142
$ws.Range("K122").Value = 2700
$ws.Range("L122").Value = 11697.4278
$ws.Range("M122").Value = -250
$ws.Range("N122").Value = -16597.4278
$ws.Range("H131").Value = 5953753
$ws.Range("J131").Value = 6061997.5
$ws.Range("L131").Value = 18185992.5
$ws.Range("N131").Value = -18196072.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 44768.75
$ws.Range("J133").Value = 44768.75
$ws.Range("L133").Value = 44768.75
$ws.Range("N133").Value = -54888.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 9419.083000000001
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 12225.444
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 12225.444
$ws.Range("M22").Value = -705
$ws.Range("N22").Value = -12815.444
$ws.Range("H27").Value = 9419.083000000001
$ws.Range("I27").Value = 1000
$ws.Range("J27").Value = 12225.444
$ws.Range("K27").Value = 1000
$ws.Range("L27").Value = 12225.444
$ws.Range("M27").Value = -893
$ws.Range("N27").Value = -12439.444
$ws.Range("H61").Value = 1213.0358
$ws.Range("I61").Value = 927.7083
$ws.Range("K61").Value = 927.7083
$ws.Range("M61").Value = -725.7083
$ws.Range("H68").Value = 2286
$ws.Range("I68").Value = 1750.5
$ws.Range("K68").Value = 1750.5
$ws.Range("M68").Value = -1001.5
$ws.Range("H71").Value = 2286
$ws.Range("I71").Value = 1750.5
$ws.Range("K71").Value = 8752.5
$ws.Range("M71").Value = -5008.5
$ws.Range("H113").Value = 1213.0358
$ws.Range("I113").Value = 927.7083
$ws.Range("K113").Value = 927.7083
$ws.Range("M113").Value = 1242.2917
$ws.Range("H136").Value = 4986.5
$ws.Range("I136").Value = 3142.4583
$ws.Range("J136").Value = 9412.200000000001
$ws.Range("K136").Value = 9427.374899999999
$ws.Range("L136").Value = 28236.6
$ws.Range("M136").Value = -6877.374899999999
$ws.Range("N136").Value = -33336.60000000001
$ws.Range("H138").Value = 21619.334
$ws.Range("J138").Value = 21619.334
$ws.Range("L138").Value = 21619.334
$ws.Range("N138").Value = -31899.334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3001.0322
$ws.Range("I132").Value = 3086.4
$ws.Range("J132").Value = 2921
$ws.Range("K132").Value = 9259.200000000001
$ws.Range("L132").Value = 8763
$ws.Range("M132").Value = -6729.200000000001
$ws.Range("N132").Value = -13823
